$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-03-05 Tuesday"; new = "2024-03-06 Wednesday"},
    @{old = "620÷5="; new = "428÷9="},
    @{old = "341÷3="; new = "561÷6="},
    @{old = "791÷5="; new = "788÷5="},
    @{old = "497÷9="; new = "687÷9="},
    @{old = "156÷3="; new = "879÷6="},
    @{old = "731÷9="; new = "355÷5="},
    @{old = "414÷4="; new = "593÷3="},
    @{old = "618÷4="; new = "601÷7="},
    @{old = "130÷3="; new = "926÷7="},
    @{old = "409÷9="; new = "621÷6="},
    @{old = "131÷7="; new = "225÷2="},
    @{old = "441÷6="; new = "250÷7="},
    @{old = "436÷8="; new = "387÷7="},
    @{old = "471÷4="; new = "107÷2="},
    @{old = "745÷4="; new = "336÷3="},
    @{old = "919÷9="; new = "517÷5="},
    @{old = "221÷7="; new = "337÷4="},
    @{old = "696÷9="; new = "632÷4="},
    @{old = "687÷8="; new = "125÷3="},
    @{old = "513÷9="; new = "487÷5="},
    @{old = "721÷6="; new = "769÷5="},
    @{old = "490÷9="; new = "993÷9="},
    @{old = "881÷2="; new = "918÷6="},
    @{old = "203÷8="; new = "689÷8="},
    @{old = "839÷3="; new = "576÷6="}
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
